$wb = $excel.ActiveWorkbook
$ws = $wb.Worksheets.Item("Sheet1")

# Insert a new row at row 7 (shifts everything below down by one),
# and populate it with the new FxE "output" / "configuration_fxe" record.
[void]$ws.Rows("7:7").Insert()

$ws.Range("A7").Value = "CHE"
$ws.Range("B7").Value = "ext_hydro"
$ws.Range("C7").Value = "output"
$ws.Range("D7").Value = "configuration_fxe"
$ws.Range("F7").Value = "hydro"
$ws.Range("G7").Value = 1

# The autofilter range and the hidden _FilterDatabase defined name both
# covered one extra row beyond the data (A5:L572); grow them by one row
# to A5:L573 to keep pace with the inserted row.
if ($ws.AutoFilterMode) {
    $ws.AutoFilterMode = $false
}
[void]$ws.Range("A5:L573").AutoFilter()

foreach ($n in $wb.Names) {
    if ($n.Name -eq "Sheet1!_FilterDatabase") {
        $n.RefersTo = "=Sheet1!`$A`$5:`$L`$573"
    }
}

# Update the selection to match the new active cell reported in the diff.
[void]$ws.Range("E7").Select()
